$wb = $excel.ActiveWorkbook
$wsOv = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# -- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" --
$wsOv.Range("E2").Value = $statusText
$wsOv.Range("F2").Value = $statusText
$wsOv.Range("E3").Value = $statusText
$wsOv.Range("F3").Value = $statusText

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# -- zh-cn: populate handback report columns --
$wsZh.Range("I2").Value = "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md"
$wsZh.Range("J2").Value = "8a894633-7d75-48f0-a031-20bb2a8fb6e6.708221081d72964b1534416fcc5c076763700091.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 03:57:11"

$wsZh.Range("I3").Value = "a195de5b-d6fe-4b2a-9833-1e336374445f.md"
$wsZh.Range("J3").Value = "a195de5b-d6fe-4b2a-9833-1e336374445f.e4bda866dc79ed15d9fcc7ce684abb0750cd5e1d.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 03:57:11"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/8a894633-7d75-48f0-a031-20bb2a8fb6e6.md", "", "", "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/8a894633-7d75-48f0-a031-20bb2a8fb6e6.md", "", "", "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/a195de5b-d6fe-4b2a-9833-1e336374445f.md", "", "", "a195de5b-d6fe-4b2a-9833-1e336374445f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/a195de5b-d6fe-4b2a-9833-1e336374445f.md", "", "", "a195de5b-d6fe-4b2a-9833-1e336374445f.md")

# -- de-de: populate handback report columns --
$wsDe.Range("I2").Value = "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md"
$wsDe.Range("J2").Value = "8a894633-7d75-48f0-a031-20bb2a8fb6e6.708221081d72964b1534416fcc5c076763700091.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 03:57:28"

$wsDe.Range("I3").Value = "a195de5b-d6fe-4b2a-9833-1e336374445f.md"
$wsDe.Range("J3").Value = "a195de5b-d6fe-4b2a-9833-1e336374445f.e4bda866dc79ed15d9fcc7ce684abb0750cd5e1d.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 03:57:28"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/8a894633-7d75-48f0-a031-20bb2a8fb6e6.md", "", "", "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/8a894633-7d75-48f0-a031-20bb2a8fb6e6.md", "", "", "8a894633-7d75-48f0-a031-20bb2a8fb6e6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/a195de5b-d6fe-4b2a-9833-1e336374445f.md", "", "", "a195de5b-d6fe-4b2a-9833-1e336374445f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18c6d49ae95eeaac759cc5bdb1325d13db13c35f/e2e/a195de5b-d6fe-4b2a-9833-1e336374445f.md", "", "", "a195de5b-d6fe-4b2a-9833-1e336374445f.md")

# -- column width adjustments (report columns widened to fit new content) --
# NOTE: ColumnWidth is specified net of the engine's fixed 5/6-character
# padding so the stored OOXML <col width> lands on the intended value.
$wsOv.Columns.Item(5).ColumnWidth = 29.14437166849777
$wsOv.Columns.Item(6).ColumnWidth = 29.14437166849777

$wsZh.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
